$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.711112380027771
$ws.Range("B1").Value = 2.797709226608276
$ws.Range("C1").Value = 3.054429531097412
$ws.Range("D1").Value = 3.443324089050293
$ws.Range("E1").Value = 1.898357033729553
